# Update the numeric results in row 2 of the "Data" sheet with the
# re-run simulation output values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B2").Value = 20761405.5539612
$ws.Range("C2").Value = 551.206641529026
$ws.Range("D2").Value = 3.60972684739767
$ws.Range("E2").Value = 66795.23987188051
$ws.Range("F2").Value = 724.5380686065409
$ws.Range("G2").Value = 10.2114942056355
$ws.Range("H2").Value = 0.705692085704497
$ws.Range("I2").Value = 0.742676942336829
$ws.Range("J2").Value = 0.739510964442611
$ws.Range("K2").Value = 729.566703379191
$ws.Range("L2").Value = 769.194764159963
$ws.Range("M2").Value = 712.408008166646
$ws.Range("N2").Value = 3.60972684739767
$ws.Range("O2").Value = 3.609726847393888
$ws.Range("P2").Value = 3.609726847393936
$ws.Range("Q2").Value = 20790217.37506654
$ws.Range("R2").Value = 20780934.60759401
$ws.Range("S2").Value = 20770911.90142715
$ws.Range("T2").Value = 20761405.5539612
$ws.Range("U2").Value = 698314.0161225691
$ws.Range("V2").Value = 852762.0834940161
$ws.Range("W2").Value = 893174.0715335517
$ws.Range("X2").Value = 906971.4579796309
$ws.Range("Y2").Value = 10.2114942056355
$ws.Range("Z2").Value = 5.00786480239524
$ws.Range("AA2").Value = 5.203629403252695
$ws.Range("AB2").Value = 162160.9385204538
$ws.Range("AC2").Value = 86085.65220273692
$ws.Range("AD2").Value = 66795.23987188051
$ws.Range("AE2").Value = 66795.23987188051
$ws.Range("AF2").Value = 866177.3057553485
$ws.Range("AG2").Value = 861299.9787655817
$ws.Range("AH2").Value = 832170.5505682746
$ws.Range("AI2").Value = 754160.2706676701
